# Update "want to go" counts (column F) across the four sheets to match
# the newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 246
$ws.Range("F5").Value  = 9432
$ws.Range("F6").Value  = 9432
$ws.Range("F7").Value  = 605
$ws.Range("F10").Value = 282
$ws.Range("F14").Value = 439
$ws.Range("F15").Value = 12061
$ws.Range("F16").Value = 12061
$ws.Range("F20").Value = 33
$ws.Range("F22").Value = 154
$ws.Range("F24").Value = 239
$ws.Range("F26").Value = 25
$ws.Range("F33").Value = 69
$ws.Range("F37").Value = 1002
$ws.Range("F38").Value = 4195
$ws.Range("F39").Value = 3634
$ws.Range("F40").Value = 540
$ws.Range("F41").Value = 2622
$ws.Range("F44").Value = 195
$ws.Range("F46").Value = 422
$ws.Range("F47").Value = 524
$ws.Range("F49").Value = 224

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 48
$ws.Range("F17").Value = 14

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 50

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value  = 246
$ws.Range("F9").Value  = 9432
$ws.Range("F10").Value = 605
$ws.Range("F13").Value = 282
$ws.Range("F17").Value = 12061
$ws.Range("F18").Value = 12061
$ws.Range("F20").Value = 33
$ws.Range("F23").Value = 50
$ws.Range("F24").Value = 154
$ws.Range("F34").Value = 69
$ws.Range("F39").Value = 1002
$ws.Range("F42").Value = 3634
$ws.Range("F46").Value = 195
$ws.Range("F47").Value = 422
$ws.Range("F49").Value = 524
$ws.Range("F51").Value = 224
